# The underlying OOXML commit swaps the two embedded DrawingML themes:
#   ppt/theme/theme1.xml (used by the slide master / all 18 slides) goes
#   from the "Integral" colour scheme to the stock "Office Theme" colour
#   scheme, while ppt/theme/theme2.xml (used only by the notes master)
#   picks up the colours the slide theme used to have.
#
# The fmtScheme (fills/lines/effects) and fontScheme (Arial everywhere)
# blocks are already byte-identical between the two themes, so the only
# observable difference is the 12-slot colour scheme (clrScheme). We
# drive that through the Theme's ThemeColorScheme collection on the
# slide master, which is the piece of the OOXML this host's object
# model actually lets us rewrite.

$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$cs = $sm.Theme.ThemeColorScheme

# Index -> target "Office Theme" RGB (COM long, 0x00BBGGRR order)
# 1 dk1      000000 -> 0
# 2 lt1      FFFFFF -> 16777215
# 3 dk2      44546A -> 6968388
# 4 lt2      E7E6E6 -> 15132391
# 5 accent1  5B9BD5 -> 13998939
# 6 accent2  ED7D31 -> 3243501
# 7 accent3  A5A5A5 -> 10855845
# 8 accent4  FFC000 -> 49407
# 9 accent5  4472C4 -> 12874308
# 10 accent6 70AD47 -> 4697456
# 11 hlink    0563C1 -> 12673797
# 12 folHlink 954F72 -> 7491477

$cs.Item(1).RGB = 0
$cs.Item(2).RGB = 16777215
$cs.Item(3).RGB = 6968388
$cs.Item(4).RGB = 15132391
$cs.Item(5).RGB = 13998939
$cs.Item(6).RGB = 3243501
$cs.Item(7).RGB = 10855845
$cs.Item(8).RGB = 49407
$cs.Item(9).RGB = 12874308
$cs.Item(10).RGB = 4697456
$cs.Item(11).RGB = 12673797
$cs.Item(12).RGB = 7491477
